# daily auto push: 2026-02-06 07:10 UTC
# Prepend a new daily-stats row (2026/02/06, 金, 14, 71) above the existing
# "2026/12/29" block on Sheet1, shifting the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 764; everything that was at row 764
# downward (through the former last row 805) shifts down to 765..806.
$ws.Rows.Item(764).Insert()

# Column A holds a plain text date string ("YYYY/MM/DD"), not a real Excel
# date. Force text formatting first so the assignment isn't auto-parsed
# into a date serial, then drop the custom number format again so the
# cell ends up with no special style, matching its siblings.
$ws.Cells.Item(764, 1).NumberFormat = "@"
$ws.Cells.Item(764, 1).Value = "2026/02/06"
$ws.Cells.Item(764, 1).ClearFormats()

$ws.Cells.Item(764, 2).Value = "金"
$ws.Cells.Item(764, 3).Value = 14
$ws.Cells.Item(764, 4).Value = 71
